# Fruta / hortaliza, semanal
# Insert a new data row at row 125 (Vega Modelo de Temuco - Pomelo, Start Ruby / Primera),
# which pushes the existing rows 125-192 down to 126-193, and fill in the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 125 - this shifts rows 125:192 down to 126:193
# and copies formatting (incl. the date numFmt style) from the row above, matching row 124/126.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new record's data.
$ws.Range("A125").Value = 10
$ws.Range("B125").Value = "Vega Modelo de Temuco"
$ws.Range("C125").Value = "La Araucanía"
$ws.Range("D125").Value = 44596
$ws.Range("E125").Value = 9
$ws.Range("F125").Value = "Fruta"
$ws.Range("G125").Value = 100102
$ws.Range("H125").Value = "Cítricos"
$ws.Range("I125").Value = 100102006
$ws.Range("J125").Value = "Pomelo"
$ws.Range("K125").Value = "Start Ruby"
$ws.Range("L125").Value = "Primera"
$ws.Range("M125").Value = 35
$ws.Range("N125").Value = 17000
$ws.Range("O125").Value = 17000
$ws.Range("P125").Value = 17000
$ws.Range("Q125").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R125").Value = "Región de O'Higgins"
$ws.Range("S125").Value = 1133
$ws.Range("T125").Value = 15
